$wb = $excel.ActiveWorkbook

# --- Sheet "Top 10 players goal 90" : rerun of xdribble model -> updated R-Proposed (D) values ---
$ws2 = $wb.Worksheets.Item("Top 10 players goal 90")
$ws2.Range("D2").Value = 77
$ws2.Range("D3").Value = 102
$ws2.Range("D4").Value = 108
$ws2.Range("D5").Value = 49
$ws2.Range("D6").Value = 89
$ws2.Range("D7").Value = 58
$ws2.Range("D8").Value = 103
$ws2.Range("D9").Value = 43
$ws2.Range("D10").Value = 54
$ws2.Range("D11").Value = 20
[void]$ws2.Range("D12").Select()

# --- Sheet "Top 10 players assist 90" : rerun of xdribble model -> updated R-Proposed (D) values ---
$ws3 = $wb.Worksheets.Item("Top 10 players assist 90")
$ws3.Range("D2").Value = 6
$ws3.Range("D3").Value = 35
$ws3.Range("D4").Value = 49
$ws3.Range("D5").Value = 29
$ws3.Range("D7").Value = 85
$ws3.Range("D8").Value = 27
$ws3.Range("D9").Value = 58
$ws3.Range("D10").Value = 20
$ws3.Range("D11").Value = 59
[void]$ws3.Range("D12").Select()

# --- Sheet "Top 10 players goal assist 90" : rerun of xdribble model -> updated R-Proposed (D) values ---
$ws4 = $wb.Worksheets.Item("Top 10 players goal assist 90")
$ws4.Range("D2").Value = 6
$ws4.Range("D3").Value = 49
$ws4.Range("D4").Value = 77
$ws4.Range("D5").Value = 58
$ws4.Range("D6").Value = 102
$ws4.Range("D7").Value = 108
$ws4.Range("D8").Value = 35
$ws4.Range("D9").Value = 89
$ws4.Range("D10").Value = 103
$ws4.Range("D11").Value = 29
[void]$ws4.Range("D12").Select()
